# Reduced the starting value for the questConditions file from 80 to 60,
# as some participants are not thresholding for the quiet task.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# startVal (B2) goes from 80 -> 60; dependent formula in C2 (=0.4*B2)
# recalculates automatically to 24.
$ws.Range("B2").Value = 60

# Leave the selection where the author left it when the file was saved.
$ws.Range("B3").Select()
